# Add two new daily NRFI-tracker sheets: 07-26-24 and 07-27-24.
# Each sheet gets a bold/bordered "Games" / "Score" header row followed by
# that day's games sorted by score (high to low), matching the layout used
# by every other dated sheet already in this workbook.

$wb = $excel.ActiveWorkbook

function Add-NrfiSheet {
    param(
        [string]$SheetName,
        [object[]]$Rows
    )

    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $SheetName

    $ws.Range("A1").Value = "Games"
    $ws.Range("B1").Value = "Score"

    $headerRange = $ws.Range("A1:B1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
    $headerRange.Borders.LineStyle = 1

    $r = 2
    foreach ($pair in $Rows) {
        $ws.Cells.Item($r, 1).Value = $pair[0]
        $ws.Cells.Item($r, 2).Value = $pair[1]
        $r = $r + 1
    }

    $null = $ws.Range("A1").Select()
}

$data21 = @(
    ,@("('CWS', 'SEA')", 0.732)
    ,@("('CHC', 'KC')", 0.731)
    ,@("('HOU', 'LAD')", 0.726)
    ,@("('CLE', 'PHI')", 0.725)
    ,@("('AZ', 'PIT')", 0.723)
    ,@("('LAA', 'OAK')", 0.719)
    ,@("('CIN', 'TB')", 0.718)
    ,@("('DET', 'MIN')", 0.703)
    ,@("('COL', 'SF')", 0.697)
    ,@("('STL', 'WSH')", 0.673)
    ,@("('ATL', 'NYM')", 0.67)
    ,@("('TEX', 'TOR')", 0.667)
    ,@("('MIA', 'MIL')", 0.5639999999999999)
    ,@("('BAL', 'SD')", 0.556)
    ,@("('BOS', 'NYY')", 0.045)
)

$data22 = @(
    ,@("('CWS', 'SEA')", 0.748)
    ,@("('BAL', 'SD')", 0.734)
    ,@("('AZ', 'PIT')", 0.73)
    ,@("('CHC', 'KC')", 0.728)
    ,@("('CIN', 'TB')", 0.713)
    ,@("('STL', 'WSH')", 0.7)
    ,@("('HOU', 'LAD')", 0.699)
    ,@("('MIA', 'MIL')", 0.697)
    ,@("('BOS', 'NYY')", 0.6820000000000001)
    ,@("('CLE', 'PHI')", 0.674)
    ,@("('COL', 'SF')", 0.575)
    ,@("('COL2', 'SF2')", 0.575)
    ,@("('TEX', 'TOR')", 0.548)
    ,@("('DET', 'MIN')", 0.479)
    ,@("('ATL', 'NYM')", 0.281)
    ,@("('LAA', 'OAK')", 0.245)
)

Add-NrfiSheet "07-26-24" $data21
Add-NrfiSheet "07-27-24" $data22

Write-Host "Total sheets:" $wb.Worksheets.Count
